$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "35.461.07"
$ws.Range("E2").Value = "  +0.68%  "

# Row 3
$ws.Range("D3").Value = "1.921.45"
$ws.Range("E3").Value = "  +1.62%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.731"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +11.96%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "254.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.69%  "

# Row 7
$ws.Range("E7").Value = "  +0.09%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.84"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.83%  "

# Row 9
$ws.Range("E9").Value = "  +2.99%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.21%  "

# Row 11
$ws.Range("E11").Value = "  +5.88%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0996"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.11%  "

# Row 13
$ws.Range("D13").Value = "2.202.09"
$ws.Range("E13").Value = "  +1.73%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.38%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.722"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.62%  "

# Row 16
$ws.Range("D16").Value = "1.930.03"
$ws.Range("E16").Value = "  +1.27%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.76%  "

# Row 18
$ws.Range("D18").Value = "35.472.98"
$ws.Range("E18").Value = "  +0.77%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.62%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0841"
$ws.Range("E20").Value = "  +3.65%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "243.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.46%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.85%  "

# Row 23
$ws.Range("E23").Value = "  +8.32%  "

# Row 24
$ws.Range("E24").Value = "  +0.05%  "

# Row 25
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.29%  "

# Row 26
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.46%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.44%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.15%  "

# Row 29
$ws.Range("E29").Value = "  +7.09%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.97%  "

# Row 31
$ws.Range("D31").Value = "4.126.04"
$ws.Range("E31").Value = "  +19.38%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.71%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +15.28%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +23.71%  "

# Row 35
$ws.Range("E35").Value = "  +4.36%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.45%  "

# Row 37
$ws.Range("E37").Value = "  +0.06%  "

# Row 38
$ws.Range("E38").Value = "  -2.24%  "

# Row 39
$ws.Range("E39").Value = "  +0.87%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.50"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.02%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.27%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.52%  "

# Row 43
$ws.Range("E43").Value = "  +1.18%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0656"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.99%  "

# Row 45
$ws.Range("D45").Value = "1.347.85"
$ws.Range("E45").Value = "  +0.96%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.84%  "

# Row 47
$ws.Range("E47").Value = "  +1.27%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.25%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.54%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.18%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.54%  "
